$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: the write order below matters - it controls the order in which new
# strings are appended to the shared-string table, which must match the
# target file's string indices.

# --- Row 13: hours updated ---
$ws.Range("E13").Value = 3

# --- Row 17 (new date row) ---
$ws.Range("A17").Value = "2019-06-09"
$ws.Range("A17").NumberFormat = "d-mmm"
$ws.Range("B17").Value = "PlayerCollision.cs"
$ws.Range("E17").Value = 2

# --- Row 18 (new date row) ---
$ws.Range("A18").Value = "2019-06-09"
$ws.Range("A18").NumberFormat = "d-mmm"
$ws.Range("B18").Value = "Pelaajan kulkualustan fysiikka"
$ws.Range("E18").Value = 1

# --- Row 11: "Pelidemo" -> "Peli" ---
$ws.Range("B11").Value = "Peli"

# --- Row 19 (new date row) ---
$ws.Range("A19").Value = "2019-06-12"
$ws.Range("A19").NumberFormat = "d-mmm"
$ws.Range("B19").Value = "score.cs ja level_complete.cs"
$ws.Range("E19").Value = 4

# --- Row 20 (new date row) ---
$ws.Range("A20").Value = "2019-06-12"
$ws.Range("A20").NumberFormat = "d-mmm"
$ws.Range("B20").Value = "levelcomplete-animaatio"
$ws.Range("E20").Value = 1

# --- Row 21 (new date row) ---
$ws.Range("A21").Value = "2019-06-14"
$ws.Range("A21").NumberFormat = "d-mmm"
$ws.Range("B21").Value = "EndTrigger.cs"
$ws.Range("E21").Value = 2

# --- Row 22 (new date row) ---
$ws.Range("A22").Value = "2019-06-18"
$ws.Range("A22").NumberFormat = "d-mmm"
$ws.Range("B22").Value = "Taso kaksi terrain"
$ws.Range("E22").Value = 4

# --- Row 24 ---
$ws.Range("B24").Value = "Peli #2 (Varasuunnitelma)"

# --- Row 26 (new date row) ---
$ws.Range("A26").Value = "2019-07-02"
$ws.Range("A26").NumberFormat = "d-mmm"
$ws.Range("B26").Value = "Varasuunnitelman laatiminen"
$ws.Range("E26").Value = 4

# --- Row 27 (new date row) ---
$ws.Range("A27").Value = "2019-07-04"
$ws.Range("A27").NumberFormat = "d-mmm"
$ws.Range("B27").Value = "Grafiikan teko"
$ws.Range("E27").Value = 10

# --- Row 28 (new date row) ---
$ws.Range("A28").Value = "2019-07-07"
$ws.Range("A28").NumberFormat = "d-mmm"
$ws.Range("B28").Value = "Perusmekaniikan teko"
$ws.Range("E28").Value = 6

# --- Row 29 (new date row) ---
$ws.Range("A29").Value = "2019-07-09"
$ws.Range("A29").NumberFormat = "d-mmm"
$ws.Range("B29").Value = "Tutoriaalitason terrain"
$ws.Range("E29").Value = 6

# --- Row 30 (new date row) ---
$ws.Range("A30").Value = "2019-07-09"
$ws.Range("A30").NumberFormat = "d-mmm"
$ws.Range("B30").Value = "Tekstit ja niiden ohjeistusgrafiikka"
$ws.Range("E30").Value = 2

# --- Row 31 (new date row) ---
$ws.Range("A31").Value = "2019-07-12"
$ws.Range("A31").NumberFormat = "d-mmm"
$ws.Range("B31").Value = "Pelaajahahmo (animointi yms)"
$ws.Range("E31").Value = 8

# --- Row 32 (new date row) ---
$ws.Range("A32").Value = "2019-07-13"
$ws.Range("A32").NumberFormat = "d-mmm"
$ws.Range("B32").Value = "Vihollishahmo (animointi yms)"
$ws.Range("E32").Value = 7

# --- Row 33 (new date row) ---
$ws.Range("A33").Value = "2019-07-16"
$ws.Range("A33").NumberFormat = "d-mmm"
$ws.Range("B33").Value = "LoadingScreen"
$ws.Range("E33").Value = 1

# --- Row 34 (new date row) ---
$ws.Range("A34").Value = "2019-07-16"
$ws.Range("A34").NumberFormat = "d-mmm"
$ws.Range("B34").Value = "Taso 1 Terrain "
$ws.Range("E34").Value = 6

# --- Row 35 (new date row) ---
$ws.Range("A35").Value = "2019-07-22"
$ws.Range("A35").NumberFormat = "d-mmm"
$ws.Range("B35").Value = "Musiikin suunnittelu"
$ws.Range("E35").Value = 2

# --- Row 36: sum formula ---
$ws.Range("E36").Formula = "=E26+E27+E28+E29+E30+E31+E32+E33+F34+E34+E35"

# --- Row 37 ---
$ws.Range("A37").Value = "MUSIIKKI VIELÄ KESKEN!!!"
$ws.Range("D37").Value = "Musiikille suunniteltu aika: 17"

# --- Row 39 ---
$ws.Range("A39").Value = "ÄÄNIEFEKTIT KESKEN!!!"
$ws.Range("D39").Value = "Ääniefekteille suunniteltu aika 6"

# --- selection / view state ---
$ws.Range("E39").Select()
